# Edit: Thu, Jul 30, 2020 5:04:59 AM
#
# 1) Slide 6's table switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 3" table style.
# 2) The presentation's theme (theme1.xml, applied through the Slide
#    Master) is changed from the "Integral" palette to the standard
#    "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{7CEEB51D-5008-44A0-BE38-76C157B25C54}")

# --- 2) Swap the theme colour palette for "Office Theme" -----------------
$colors = $p.Slides.Item(1).ThemeColorScheme

$colors.Item(1).RGB  = 0x000000   # dk1
$colors.Item(2).RGB  = 0xFFFFFF   # lt1
$colors.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink -> 954F72
